# Replace the lattice-multiplication exercises in the 5x3 table with a
# new set of problems, rewriting each cell's 5-line content
# ("AA x BB" / "  C    D" / "  ----" / "E|    |" / "F|    |") in place.
#
# Each cell is a single run containing five <w:t> text nodes separated by
# <w:br/> line breaks. Setting Cell.Range.Text to a string with embedded
# vertical-tab characters (chr 11) reproduces that same t/br/t/br/... run
# layout, while Word preserves the existing run formatting (sz=32).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11

function Set-LatticeCell($table, $row, $col, $top, $side, $r1, $r2) {
    $table.Cell($row, $col).Range.Text = $top + $br + $side + $br + '  ----' + $br + $r1 + $br + $r2
}

Set-LatticeCell $t 1 1 '21 x 79' '  7    9' '2|    |' '1|    |'
Set-LatticeCell $t 1 2 '46 x 30' '  3    0' '4|    |' '6|    |'
Set-LatticeCell $t 1 3 '35 x 64' '  6    4' '3|    |' '5|    |'

Set-LatticeCell $t 2 1 '10 x 57' '  5    7' '1|    |' '0|    |'
Set-LatticeCell $t 2 2 '74 x 47' '  4    7' '7|    |' '4|    |'
Set-LatticeCell $t 2 3 '23 x 62' '  6    2' '2|    |' '3|    |'

Set-LatticeCell $t 3 1 '94 x 18' '  1    8' '9|    |' '4|    |'
Set-LatticeCell $t 3 2 '92 x 25' '  2    5' '9|    |' '2|    |'
Set-LatticeCell $t 3 3 '99 x 26' '  2    6' '9|    |' '9|    |'

Set-LatticeCell $t 4 1 '84 x 48' '  4    8' '8|    |' '4|    |'
Set-LatticeCell $t 4 2 '33 x 58' '  5    8' '3|    |' '3|    |'
Set-LatticeCell $t 4 3 '11 x 72' '  7    2' '1|    |' '1|    |'

Set-LatticeCell $t 5 1 '50 x 25' '  2    5' '5|    |' '0|    |'
Set-LatticeCell $t 5 2 '62 x 14' '  1    4' '6|    |' '2|    |'
Set-LatticeCell $t 5 3 '83 x 21' '  2    1' '8|    |' '3|    |'

Write-Output "Updated $($t.Rows.Count)x$($t.Columns.Count) lattice table."
